# Daily attendance processing - 2026-01-01 22:32:23
# Swap the order of "Recorded By" contributors in column G:
#   "System, dnasr281@gmail.com"  ->  "dnasr281@gmail.com, System"
# Applies to every row in the used range whose column-G text matches the
# old ordering; all other columns/cells are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldText = "System, dnasr281@gmail.com"
$newText = "dnasr281@gmail.com, System"

$lastRow = $ws.UsedRange.Rows.Count
$changed = 0

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    if ($cell.Text -eq $oldText) {
        $cell.Value = $newText
        $changed = $changed + 1
    }
}

Write-Output "Updated $changed cell(s) in column G"
